$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to be treated as text while we write the updated price/volume
# strings below (several of them are plain decimals like "1.002" that Excel
# would otherwise auto-convert to numbers), then restore the original
# (default/"Normal") cell style so no formatting is left behind.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "23.424.16"
$ws.Range("E2").Value = "  -0.40%  "

# Row 3
$ws.Range("D3").Value = "1.629.65"
$ws.Range("E3").Value = "  -0.67%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("D6").Value = "304.26"
$ws.Range("E6").Value = "  -1.10%  "

# Row 7
$ws.Range("D7").Value = "0.3780"
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("D8").Value = "0.3652"
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "51.44"
$ws.Range("E9").Value = "  -1.65%  "

# Row 10
$ws.Range("D10").Value = "0.08232"
$ws.Range("E10").Value = "  +0.63%  "

# Row 11
$ws.Range("D11").Value = "1.229"
$ws.Range("E11").Value = "  -3.33%  "

# Row 12
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13
$ws.Range("D13").Value = "22.39"
$ws.Range("E13").Value = "  -2.95%  "

# Row 14
$ws.Range("D14").Value = "6.556"
$ws.Range("E14").Value = "  -1.47%  "

# Row 15
$ws.Range("D15").Value = "0.00001252"
$ws.Range("E15").Value = "  -2.32%  "

# Row 16
$ws.Range("D16").Value = "7.319"
$ws.Range("E16").Value = "  -1.06%  "

# Row 17
$ws.Range("D17").Value = "1.630.75"
$ws.Range("E17").Value = "  -0.33%  "

# Row 18
$ws.Range("D18").Value = "94.10"
$ws.Range("E18").Value = "  -0.82%  "

# Row 19
$ws.Range("D19").Value = "0.06980"
$ws.Range("E19").Value = "  +0.12%  "

# Row 20
$ws.Range("D20").Value = "17.73"
$ws.Range("E20").Value = "  -2.88%  "

# Row 21
$ws.Range("D21").Value = "6.489"
$ws.Range("E21").Value = "  -1.07%  "

# Row 22
$ws.Range("E22").Value = "  +0.11%  "

# Row 23
$ws.Range("D23").Value = "12.72"
$ws.Range("E23").Value = "  -0.96%  "

# Row 24
$ws.Range("D24").Value = "23.433.93"
$ws.Range("E24").Value = "  -0.38%  "

# Row 25
$ws.Range("E25").Value = "  +3.17%  "

# Row 26
$ws.Range("D26").Value = "2.468"
$ws.Range("E26").Value = "  +2.02%  "

# Row 27
$ws.Range("D27").Value = "21.38"
$ws.Range("E27").Value = "  +0.38%  "

# Row 28
$ws.Range("D28").Value = "150.07"
$ws.Range("E28").Value = "  -0.91%  "

# Row 29
$ws.Range("D29").Value = "5.304"
$ws.Range("E29").Value = "  -1.04%  "

# Row 30
$ws.Range("D30").Value = "134.18"
$ws.Range("E30").Value = "  -1.17%  "

# Row 31
$ws.Range("D31").Value = "1.809.73"
$ws.Range("E31").Value = "  -0.46%  "

# Row 32
$ws.Range("D32").Value = "2.267"
$ws.Range("E32").Value = "  -3.81%  "

# Row 33
$ws.Range("D33").Value = "6.822"
$ws.Range("E33").Value = "  +0.38%  "

# Row 34
$ws.Range("D34").Value = "1.018"
$ws.Range("E34").Value = "  +5.31%  "

# Row 35
$ws.Range("D35").Value = "10.91"
$ws.Range("E35").Value = "  +5.22%  "

# Row 36
$ws.Range("D36").Value = "0.02792"
$ws.Range("E36").Value = "  -1.44%  "

# Row 37
$ws.Range("D37").Value = "0.2524"
$ws.Range("E37").Value = "  -0.62%  "

# Row 38
$ws.Range("D38").Value = "0.08746"
$ws.Range("E38").Value = "  -1.35%  "

# Row 39
$ws.Range("D39").Value = "0.07121"
$ws.Range("E39").Value = "  -3.42%  "

# Row 40
$ws.Range("D40").Value = "6.042"
$ws.Range("E40").Value = "  -2.42%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.353"
$ws.Range("E41").Value = "  -2.46%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.7046"
$ws.Range("E42").Value = "  -1.13%  "

# Row 43
$ws.Range("D43").Value = "16.26"
$ws.Range("E43").Value = "  +0.29%  "

# Row 44
$ws.Range("D44").Value = "12.21"
$ws.Range("E44").Value = "  -2.67%  "

# Row 45
$ws.Range("D45").Value = "0.6552"
$ws.Range("E45").Value = "  -0.24%  "

# Row 46
$ws.Range("D46").Value = "2.320"
$ws.Range("E46").Value = "  -1.12%  "

# Row 47
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("D48").Value = "3.980"
$ws.Range("E48").Value = "  -1.34%  "

# Row 49
$ws.Range("D49").Value = "0.08008"
$ws.Range("E49").Value = "  +0.38%  "

# Row 50
$ws.Range("D50").Value = "1.197"
$ws.Range("E50").Value = "  -1.28%  "

# Row 51
$ws.Range("D51").Value = "126.05"
$ws.Range("E51").Value = "  -2.68%  "

# Restore the default style on the Price/Volume columns
$priceVolRange.Style = "Normal"

